{"js": "// Update the date line and the 25 two-digit multiplication problems to\n// the values produced for the next day's worksheet.\nconst replacements = [\n  [\"2025-04-01 Tuesday\", \"2025-04-02 Wednesday\"],\n  [\"73\u00d755=\", \"50\u00d723=\"],\n  [\"32\u00d739=\", \"65\u00d712=\"],\n  [\"97\u00d775=\", \"21\u00d725=\"],\n  [\"82\u00d711=\", \"15\u00d756=\"],\n  [\"27\u00d742=\", \"93\u00d754=\"],\n  [\"53\u00d781=\", \"44\u00d772=\"],\n  [\"37\u00d749=\", \"70\u00d755=\"],\n  [\"40\u00d726=\", \"79\u00d728=\"],\n  [\"27\u00d758=\", \"96\u00d718=\"],\n  [\"62\u00d784=\", \"71\u00d799=\"],\n  [\"28\u00d765=\", \"92\u00d719=\"],\n  [\"47\u00d765=\", \"45\u00d765=\"],\n  [\"34\u00d774=\", \"82\u00d777=\"],\n  [\"50\u00d776=\", \"14\u00d731=\"],\n  [\"56\u00d791=\", \"98\u00d791=\"],\n  [\"31\u00d754=\", \"58\u00d793=\"],\n  [\"46\u00d753=\", \"42\u00d750=\"],\n  [\"87\u00d711=\", \"25\u00d717=\"],\n  [\"38\u00d711=\", \"23\u00d787=\"],\n  [\"87\u00d735=\", \"98\u00d760=\"],\n  [\"80\u00d774=\", \"25\u00d718=\"],\n  [\"56\u00d713=\", \"33\u00d744=\"],\n  [\"28\u00d796=\", \"57\u00d749=\"],\n  [\"68\u00d786=\", \"63\u00d713=\"],\n  [\"86\u00d781=\", \"18\u00d781=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 two-digit multiplication problems to\n# the values produced for the next day's worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-01 Tuesday\", \"2025-04-02 Wednesday\"),\n    @(\"73\u00d755=\", \"50\u00d723=\"),\n    @(\"32\u00d739=\", \"65\u00d712=\"),\n    @(\"97\u00d775=\", \"21\u00d725=\"),\n    @(\"82\u00d711=\", \"15\u00d756=\"),\n    @(\"27\u00d742=\", \"93\u00d754=\"),\n    @(\"53\u00d781=\", \"44\u00d772=\"),\n    @(\"37\u00d749=\", \"70\u00d755=\"),\n    @(\"40\u00d726=\", \"79\u00d728=\"),\n    @(\"27\u00d758=\", \"96\u00d718=\"),\n    @(\"62\u00d784=\", \"71\u00d799=\"),\n    @(\"28\u00d765=\", \"92\u00d719=\"),\n    @(\"47\u00d765=\", \"45\u00d765=\"),\n    @(\"34\u00d774=\", \"82\u00d777=\"),\n    @(\"50\u00d776=\", \"14\u00d731=\"),\n    @(\"56\u00d791=\", \"98\u00d791=\"),\n    @(\"31\u00d754=\", \"58\u00d793=\"),\n    @(\"46\u00d753=\", \"42\u00d750=\"),\n    @(\"87\u00d711=\", \"25\u00d717=\"),\n    @(\"38\u00d711=\", \"23\u00d787=\"),\n    @(\"87\u00d735=\", \"98\u00d760=\"),\n    @(\"80\u00d774=\", \"25\u00d718=\"),\n    @(\"56\u00d713=\", \"33\u00d744=\"),\n    @(\"28\u00d796=\", \"57\u00d749=\"),\n    @(\"68\u00d786=\", \"63\u00d713=\"),\n    @(\"86\u00d781=\", \"18\u00d781=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
